$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.285.49"
$ws.Range("E2").Value = "  +5.06%  "
$ws.Range("D3").Value = "3.242.18"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'576.38"
$ws.Range("E5").Value = "  +2.64%  "
$ws.Range("D6").Value = "'178.61"
$ws.Range("E6").Value = "  +6.67%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.85%  "
$ws.Range("D9").Value = "3.239.41"
$ws.Range("E9").Value = "  +2.50%  "
$ws.Range("E10").Value = "  +4.48%  "
$ws.Range("D11").Value = "'6.72"
$ws.Range("E11").Value = "  +2.07%  "
$ws.Range("E12").Value = "  +4.60%  "
$ws.Range("D13").Value = "3.804.37"
$ws.Range("E13").Value = "  +2.60%  "
$ws.Range("D15").Value = "'27.89"
$ws.Range("E15").Value = "  +2.70%  "
$ws.Range("D16").Value = "67.203.52"
$ws.Range("E16").Value = "  +4.98%  "
$ws.Range("E17").Value = "  +3.23%  "
$ws.Range("D18").Value = "3.244.96"
$ws.Range("E18").Value = "  +2.65%  "
$ws.Range("E19").Value = "  +2.73%  "
$ws.Range("D20").Value = "'13.30"
$ws.Range("E20").Value = "  +3.02%  "
$ws.Range("D21").Value = "'375.14"
$ws.Range("E21").Value = "  +6.91%  "
$ws.Range("D22").Value = "'7.59"
$ws.Range("E22").Value = "  +6.25%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "'71.05"
$ws.Range("E24").Value = "  +3.83%  "
$ws.Range("E25").Value = "  +1.41%  "
$ws.Range("D26").Value = "3.386.44"
$ws.Range("E26").Value = "  +2.63%  "
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("D28").Value = "'9.90"
$ws.Range("E28").Value = "  +4.20%  "
$ws.Range("E29").Value = "  +1.85%  "
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").Value = "'1.97"
$ws.Range("E31").Value = "  +4.51%  "
$ws.Range("E32").Value = "  +2.32%  "
$ws.Range("D33").Value = "'22.51"
$ws.Range("E33").Value = "  +2.88%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("E35").Value = "  +6.27%  "
$ws.Range("E36").Value = "  +3.05%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "'161.91"
$ws.Range("E37").Value = "  +5.79%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'1.49"
$ws.Range("E38").Value = "  +3.90%  "
$ws.Range("E39").Value = "  +5.16%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'6.96"
$ws.Range("E40").Value = "  +17.24%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.86"
$ws.Range("E41").Value = "  +10.35%  "
$ws.Range("D42").Value = "'26.81"
$ws.Range("E42").Value = "  +4.50%  "
$ws.Range("D43").Value = "'2.61"
$ws.Range("E43").Value = "  +5.24%  "
$ws.Range("D44").Value = "2.758.07"
$ws.Range("E44").Value = "  +6.32%  "
$ws.Range("E45").Value = "  +5.46%  "
$ws.Range("D46").Value = "'351.66"
$ws.Range("E46").Value = "  +10.54%  "
$ws.Range("D47").Value = "'25.71"
$ws.Range("E47").Value = "  +9.10%  "
$ws.Range("D48").Value = "'40.45"
$ws.Range("E48").Value = "  +2.85%  "
$ws.Range("D49").Value = "'0.0671"
$ws.Range("E49").Value = "  +3.76%  "
$ws.Range("E50").Value = "  +4.11%  "
$ws.Range("E51").Value = "  +1.80%  "
